# Q factor run for sg_rr_84_025 2023-12-11 16-27-03.csv data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New note row (row 74) - inserted right after row 73, before the blank row 75
$ws.Range("A74").Value = "Note I mostly only committed to Git after each run that I thought had gone correctly after adjusting prominences."

# Two intermediate Q-factor attempts for sg_rr_84_025 that only got as far as the
# prominence column before being abandoned (noise was found), each annotated in
# column U.
$ws.Range("A85").Value = "sg_rr_84_025 2023-12-11 16-27-03.csv"
$ws.Range("B85").Value = 0.01
$ws.Range("C85").Value = 1000
$ws.Range("D85").Value = 5001
$ws.Range("E85").Value = 1530
$ws.Range("F85").Value = 1570
$ws.Range("G85").Value = 0.00025
$ws.Range("H85").Value = "(approx_fsr/2)/wavelength step size"
$ws.Range("I85").Value = 1.7
$ws.Range("U85").Value = "seemed to find one peak in what looked like noise so increased prominence"

$ws.Range("A86").Value = "sg_rr_84_025 2023-12-11 16-27-03.csv"
$ws.Range("B86").Value = 0.01
$ws.Range("C86").Value = 1000
$ws.Range("D86").Value = 5001
$ws.Range("E86").Value = 1530
$ws.Range("F86").Value = 1570
$ws.Range("G86").Value = 0.00026
$ws.Range("H86").Value = "(approx_fsr/2)/wavelength step size"
$ws.Range("I86").Value = 1.7
$ws.Range("U86").Value = "seemed to find one peak in what looked like noise so increased prominence"

# Final successful run with the full calculated Q-factor results
$ws.Range("A87").Value = "sg_rr_84_025 2023-12-11 16-27-03.csv"
$ws.Range("B87").Value = 0.01
$ws.Range("C87").Value = 1000
$ws.Range("D87").Value = 5001
$ws.Range("E87").Value = 1530
$ws.Range("F87").Value = 1570
$ws.Range("G87").Value = 0.003
$ws.Range("H87").Value = "(approx_fsr/2)/wavelength step size"
$ws.Range("I87").Value = 1.7
$ws.Range("J87").Value = 1.1724242424242399
$ws.Range("K87").Value = 0.0059834098769671303
$ws.Range("L87").Value = "yes"
$ws.Range("M87").Value = 0.15082855204548801
$ws.Range("N87").Value = 0.0058058346947930497
$ws.Range("O87").Value = 10861.1687372212
$ws.Range("P87").Value = 523.81498314973601
$ws.Range("Q87").Value = 181581295.59314901
$ws.Range("R87").Value = 26275592.48443
$ws.Range("S87").Value = 84
$ws.Range("T87").Value = 0.1

# Update the view to reflect the new bottom of the table
$ws.Range("A88").Select()
$excel.ActiveWindow.ScrollRow = 73
